$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Egypt / "2" / Household Products
$ws.Range("D2").Value = 0.19645
$ws.Range("G2").Value = 0.01113089937666963
$ws.Range("H2").Value = 0.01113089937666963
$ws.Range("I2").Value = -0.0002671415850400714
$ws.Range("J2").Value = -0.0002671415850400714
$ws.Range("K2").Value = -2.02
$ws.Range("L2").Value = -0.01798753339269812
$ws.Range("U2").Value = 3.03
$ws.Range("V2").Value = 0.1333040035195777
$ws.Range("W2").Value = 0.1366134654191891
$ws.Range("X2").Value = 0.09978267300828686
$ws.Range("Y2").Value = 0.03683079241090229
$ws.Range("Z2").Value = 9.245842252593448
$ws.Range("AA2").Value = 0.6331706136769198
$ws.Range("AB2").Value = 0.09800787880243006
$ws.Range("AC2").Value = 0.5351627348744897
$ws.Range("AD2").Value = 1.091
$ws.Range("AF2").Value = 1.091
$ws.Range("AG2").Value = -1.939
$ws.Range("AH2").Value = 0.04579992443642164
$ws.Range("AI2").Value = 0.08308582743126952
$ws.Range("AJ2").Value = -0.09326150738300228
$ws.Range("AK2").Value = -0.191961191961192
$ws.Range("AM2").Value = -0.13
$ws.Range("AN2").Value = 0.8590551181102362
$ws.Range("AP2").Value = -1.526771653543307
$ws.Range("AQ2").Value = 0.230769230769231

# Row 3 - company name swap: PRCL -> MOSC
$ws.Range("B3").Value = "Misr Oil & Soap (CASE:MOSC)"
$ws.Range("D3").Value = 0.335
$ws.Range("G3").Value = 0.01912225705329153
$ws.Range("H3").Value = 0.01912225705329153
$ws.Range("I3").Value = 0.01661442006269593
$ws.Range("J3").Value = 0.01661442006269593
$ws.Range("K3").Value = 1.34
$ws.Range("L3").Value = 0.01400208986415883
$ws.Range("U3").Value = 1.95
$ws.Range("V3").Value = 0.2459016393442623
$ws.Range("W3").Value = 0.5654008438818565
$ws.Range("X3").Value = 0.09849641417084223
$ws.Range("Y3").Value = 0.4669044297110143
$ws.Range("Z3").Value = 85.06666666666666
$ws.Range("AA3").Value = 1.413333333333333
$ws.Range("AB3").Value = 0.09744941919653284
$ws.Range("AC3").Value = 1.3158839141368
$ws.Range("AD3").Value = 0.185
$ws.Range("AF3").Value = 0.185
$ws.Range("AG3").Value = -1.765
$ws.Range("AH3").Value = 0.02279728897104128
$ws.Range("AI3").Value = 0.04993252361673414
$ws.Range("AJ3").Value = -0.2862935928629359
$ws.Range("AK3").Value = -1.005698005698006
$ws.Range("AM3").Value = -0.035
$ws.Range("AN3").Value = 0.09999999999999999
$ws.Range("AP3").Value = -0.954054054054054
$ws.Range("AQ3").Value = -45.42857142857142

# Row 4 - company name swap: MOSC -> PRCL
$ws.Range("B4").Value = "The General Company for Ceramic and Porcelain Products (CASE:PRCL)"
$ws.Range("D4").Value = 0.0579
$ws.Range("G4").Value = -0.03493975903614457
$ws.Range("H4").Value = -0.03493975903614457
$ws.Range("I4").Value = -0.09759036144578313
$ws.Range("J4").Value = -0.09759036144578313
$ws.Range("K4").Value = -3.36
$ws.Range("L4").Value = -0.2024096385542168
$ws.Range("U4").Value = 1.08
$ws.Range("V4").Value = 0.07297297297297298
$ws.Range("W4").Value = -0.2921739130434782
$ws.Range("X4").Value = 0.1010689318457315
$ws.Range("Y4").Value = -0.3932428448892097
$ws.Range("Z4").Value = 1.506215406950368
$ws.Range("AA4").Value = -0.1469921059794937
$ws.Range("AB4").Value = 0.09856633840832729
$ws.Range("AC4").Value = -0.245558444387821
$ws.Range("AD4").Value = 0.906
$ws.Range("AF4").Value = 0.906
$ws.Range("AG4").Value = -0.174
$ws.Range("AH4").Value = 0.05768496116133962
$ws.Range("AI4").Value = 0.09611712285168683
$ws.Range("AJ4").Value = -0.0118966224531656
$ws.Range("AK4").Value = -0.02084831056793674
$ws.Range("AM4").Value = -0.095
$ws.Range("AN4").Value = -1.562068965517242
$ws.Range("AP4").Value = 0.3000000000000001
$ws.Range("AQ4").Value = 17.05263157894737
